$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Row 12 (existing row, gains a B12 cell)
$ws.Range("B12").Value = '“Ты как здесь оказался, у нас в деревне даже 100-а жителей не наберется, все друг друга в лицо знают” – спросил у меня выживший – “ Тоже не успел на эвакуацию?”. Мы сели за стол, этот выживший налил чай и попросил мне объяснить, кто я и как я здесь оказался. Моего собеседника звали Виктор. Я ничего не понимал, я попытался ему объяснить, что я турист, который был в недельном походе от одной деревни в другую. “Слушай, рассказываю для мяса, за те семь дней, которые ты был в своём походе, случилась эпидемия, которая поразила весь мир. Люди превратились в кровожадных тварей, которые жаждут тебя сожрать. Если так можно сказать случился обычный “Голливудский” зомби апокалипсис. Все города превратились в огромные рассадники зомби, там жизни нет. Все кто успел, эвакуировались из городов в деревни и сёла, но вирус добрался и до нас, не знаю насчёт остальных стран, но Россия уже почти полностью пала, последние новости, которые я смог увидеть были в четверг, на них говорилось, что вспышки вируса были замечены во всех странах. Они объявили, что эвакуация происходит по некоторым шоссе, в том числе и тому, которое располагается недалеко от этой деревни. Но из-за того что автобусы, которые приехали в эту деревню привлекли зомби с округи, не все успели эвакуироваться и их либо сожрали, либо они смогли убежать в неизвестном направлении, как я. Думаю что если идти вдоль шоссе то мы сможем дойти до военного блокпоста, так как звук привлекает их, а значит они должны последовать за автобусами ”. Я был шокирован услышанным, но когда с чердака избушки, в которой мы были, я увидел, как зомби пожирали чей-то труп, то мне просто пришлось поверить в данную историю. Мой спаситель был мужчина возраста около 45 лет, по его словам он был ветераном Чеченской войны, и поэтому навыки, полученные им на войне, помогли ему выжить.. Он предложил мне попробовать добраться до военного блокпоста, потому что это единственный способ выжить. Мы собрали все припасы, которые есть у Виктора. Самое главное, что Виктор был охотником и в его запасе имелись блочный лук, арбалет  и боеприпасы к ним. Он дал мне на вооружение арбалет, а сам взял лук, так как сказал, что с арбалетом мне будет обращаться легче, чем с луком. Но как только мы вышли на улицу, нам встретилось трое ходячих мертвецов. Что же нам с ними делать?'

# Row 13 (new)
$ws.Range("A13").Value = 21
$ws.Range("B13").Value = 'Пока мы собирали припасы, Виктор объяснил мне как пользоваться арбалетом из-за чего убить зомби мне не составило труда, особенно под надзором Виктора. Мы пошли к выходу из деревни, как Виктор и предполагал звуки моторов эвакуационных автобусов, уманили часть зомби за собой, а остальную привлекли за собой не успевшие  эвакуироваться люди. Из-за чего нам не составило труда выйти до шоссе и пойти по нему. '

# Row 14 (new, vertical-center aligned B cell)
$ws.Range("A14").Value = 22
$ws.Range("B14").Value = 'Мы спрятались за углом здания и стали кидать камни в сторону зомби. Поначалу таким трюком зомби не приманивались, но на один из камней они отреагировали все втроём, двое из них даже побежали. Виктор ловким движением рук выстрелил из лука и убил одного из них, но паника сыграла свою роль, и мой арбалет дал осечку зомби повалил меня и укусил за шею, я кричал от боли очень громко. Виктор убил оставшихся остальных зомби, он успокоил меня тем, что по его наблюдениям, укус заражает, только если тебя укусил кровавый зомби, но было уже слишком. На мои крики сбежалось около 30 зомби, мы пытались отстреливаться, но нам не хватило времени на то, чтобы убить их всех, они повалили нас и загрызли. Смерть '
$ws.Range("B14").VerticalAlignment = -4108

# Row 15 (new)
$ws.Range("A15").Value = 23
$ws.Range("B15").Value = 'Мы решили обойти их по-тихому, что не составило труда. Мы пошли к выходу из деревни, как Виктор и предполагал звуки моторов эвакуационных автобусов, уманили часть зомби за собой, а остальную привлекли за собой не успевшие  эвакуироваться люди. Из-за чего нам не составило труда выйти до шоссе и пойти по нему. '

# Update selection to match the commit (B3 was selected when saved)
$ws.Range("B3").Select()

